# Daily attendance processing - 2025-12-14 22:26:37
# Normalize the "Recorded By" column (G) so that the "System" marker is
# listed after the human/automation accounts that recorded attendance,
# instead of always being listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -ne $null -and $current.StartsWith("System, ")) {
        $rest = $current.Substring(8)
        $cell.Value = $rest + ", System"
    }
}
